# Capstone Project - Poster.pptx update
# - "RESEARCH OBJECTIVES" heading renamed to "PROJECT OBJECTIVES"
# - Results table picture (Imagem 25) nudged/resized

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# 1. "RESEARCH OBJECTIVES" -> "PROJECT OBJECTIVES" (Text Placeholder 3, shape id 40)
$headingShape = Get-ShapeById $s 40
$headingShape.TextFrame.TextRange.Text = "PROJECT OBJECTIVES"

# 2. Move/resize the results-table picture (Imagem 25, shape id 26, creationId
#    {2F68D4EE-63B4-B304-F1DF-2719624F6DC5}) to its new position/size.
$pic = Get-ShapeById $s 26
$pic.Left = 33867632 / 12700
$pic.Top = 18232036 / 12700
$pic.Width = 8981310 / 12700
$pic.Height = 6201481 / 12700
